$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4459
$ws.Range("B4").Value = 14
$ws.Range("B5").Value = 85
$ws.Range("B6").Value = 4
